# ComposicionQuincenalGNALoteIV.xlsx - "primera Quincena" template rework.
#
# The report's data-source item properties were renamed (prefixed so the
# "Gna" composition item fields don't collide with the other composition
# blocks the same model now feeds, e.g. item.Fecha -> item.CompGnaDia,
# item.C6 -> item.CompGnaC6, item.Simbolo -> item.CompSimbolo, ...), and a
# brand-new row of merge-field placeholders (row 6) was added under the
# header row to hold the "Total Promedio PeruPetro" averages per component.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 3: header merge-field placeholders, renamed in place ---------------
$ws.Range("B3").Value  = "{{item.CompGnaDia}}"
$ws.Range("C3").Value  = "{{item.CompGnaC6}}"
$ws.Range("D3").Value  = "{{item.CompGnaC3}}"
$ws.Range("E3").Value  = "{{item.CompGnaIc4}}"
$ws.Range("F3").Value  = "{{item.CompGnaNc4}}"
$ws.Range("G3").Value  = "{{item.CompGnaNeoC5}}"
$ws.Range("H3").Value  = "{{item.CompGnaIc5}}"
$ws.Range("I3").Value  = "{{item.CompGnaNc5}}"
$ws.Range("J3").Value  = "{{item.CompGnaNitrog}}"
$ws.Range("K3").Value  = "{{item.CompGnaC1}}"
$ws.Range("L3").Value  = "{{item.CompGnaCo2}}"
$ws.Range("M3").Value  = "{{item.CompGnaC2}}"
$ws.Range("N3").Value  = "{{item.CompGnaObservacion}}"

# --- Row 6: new "Total Promedio PeruPetro" placeholder row ------------------
# Give it the same look as the row-3/row-4 header band before filling it in.
$ws.Range("C3:M3").Copy()
$ws.Range("C6:M6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("N3").Copy()
$ws.Range("N6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("C6").Value = "{{TotalPromedioPeruPetroC6}}"
$ws.Range("D6").Value = "{{TotalPromedioPeruPetroC3}}"
$ws.Range("E6").Value = "{{TotalPromedioPeruPetroIc4}}"
$ws.Range("F6").Value = "{{TotalPromedioPeruPetroNc4}}"
$ws.Range("G6").Value = "{{TotalPromedioPeruPetroNeoC5}}"
$ws.Range("H6").Value = "{{TotalPromedioPeruPetroIc5}}"
$ws.Range("I6").Value = "{{TotalPromedioPeruPetroNc5}}"
$ws.Range("J6").Value = "{{TotalPromedioPeruPetroNitrog}}"
$ws.Range("K6").Value = "{{TotalPromedioPeruPetroC1}}"
$ws.Range("L6").Value = "{{TotalPromedioPeruPetroCo2}}"
$ws.Range("M6").Value = "{{TotalPromedioPeruPetroC2}}"

# --- Row 10 ("Componente_Items") merge fields, renamed in place -------------
$ws.Range("B10").Value = "{{item.CompSimbolo}}"
$ws.Range("C10").Value = "{{item.CompDescripcion}}"
$ws.Range("E10").Value = "{{item.CompMolPorc}}"

# --- Restore the view to show the top of the sheet with the new row 6 cell
#     selected, same as the author left it.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I6").Select()
